# Data refresh: update "想去人数" (want-to-go count) figures in column F
# across all four worksheets, matching the upstream scrape regeneration
# ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 55
$ws.Range("F4").Value = 7959
$ws.Range("F6").Value = 8689
$ws.Range("F7").Value = 5832
$ws.Range("F9").Value = 2860
$ws.Range("F10").Value = 1241
$ws.Range("F14").Value = 660
$ws.Range("F15").Value = 136
$ws.Range("F16").Value = 4104
$ws.Range("F17").Value = 4104
$ws.Range("F18").Value = 88
$ws.Range("F19").Value = 76
$ws.Range("F20").Value = 84
$ws.Range("F22").Value = 186
$ws.Range("F23").Value = 42
$ws.Range("F24").Value = 6073
$ws.Range("F25").Value = 6073
$ws.Range("F26").Value = 219
$ws.Range("F27").Value = 84
$ws.Range("F28").Value = 291
$ws.Range("F29").Value = 412
$ws.Range("F30").Value = 189
$ws.Range("F31").Value = 435
$ws.Range("F32").Value = 4305
$ws.Range("F33").Value = 1582
$ws.Range("F35").Value = 1733
$ws.Range("F36").Value = 5782
$ws.Range("F37").Value = 90
$ws.Range("F39").Value = 70
$ws.Range("F40").Value = 50
$ws.Range("F41").Value = 3816
$ws.Range("F42").Value = 52
$ws.Range("F43").Value = 64
$ws.Range("F45").Value = 2372
$ws.Range("F50").Value = 251
$ws.Range("F51").Value = 21

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 177
$ws.Range("F4").Value = 21

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1393

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1393
$ws.Range("F4").Value = 7959
$ws.Range("F6").Value = 8689
$ws.Range("F7").Value = 5832
$ws.Range("F9").Value = 2860
$ws.Range("F10").Value = 1241
$ws.Range("F14").Value = 177
$ws.Range("F15").Value = 660
$ws.Range("F16").Value = 137
$ws.Range("F17").Value = 4104
$ws.Range("F18").Value = 4104
$ws.Range("F19").Value = 88
$ws.Range("F20").Value = 76
$ws.Range("F21").Value = 84
$ws.Range("F23").Value = 186
$ws.Range("F24").Value = 42
$ws.Range("F25").Value = 6073
$ws.Range("F26").Value = 6073
$ws.Range("F27").Value = 219
$ws.Range("F28").Value = 84
$ws.Range("F29").Value = 412
$ws.Range("F30").Value = 189
$ws.Range("F31").Value = 435
$ws.Range("F33").Value = 4305
$ws.Range("F34").Value = 1582
$ws.Range("F37").Value = 1733
$ws.Range("F39").Value = 5782
$ws.Range("F40").Value = 90
$ws.Range("F42").Value = 3816
$ws.Range("F43").Value = 64
$ws.Range("F47").Value = 2372
$ws.Range("F50").Value = 251
